$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E1").Value = "Type_prod"
$ws.Range("E1").Style = $ws.Range("D1").Style
